$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new "2021年" row (row 13) below the existing data (rows 2-12).

# --- Column A: year label, styled like the other year cells (A2:A12) ---
$ws.Range("A13").Value = "2021年"
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Numeric columns with known 2021 values ---
$ws.Range("B13").Value = 9.199999999999999
$ws.Range("C13").Value = 74.59999999999999
$ws.Range("E13").Value = 272.3

# --- Columns without data for 2021 (D, F, G, H, I): blank-but-present text
# cells, matching the pattern already used elsewhere in the sheet for
# missing data points. A bare quote forces Excel to store an empty text
# value instead of silently skipping the (no-op) empty assignment; then we
# paste plain formatting over it (from an untouched, unstyled cell) so the
# transient "quote prefix" flag doesn't linger and the cell stays unstyled,
# like its counterparts in other rows.
$ws.Range("Z1").Copy() | Out-Null

$ws.Range("D13").Value = "'"
$ws.Range("D13").PasteSpecial(-4122) | Out-Null

$ws.Range("F13").Value = "'"
$ws.Range("F13").PasteSpecial(-4122) | Out-Null

$ws.Range("G13").Value = "'"
$ws.Range("G13").PasteSpecial(-4122) | Out-Null

$ws.Range("H13").Value = "'"
$ws.Range("H13").PasteSpecial(-4122) | Out-Null

$ws.Range("I13").Value = "'"
$ws.Range("I13").PasteSpecial(-4122) | Out-Null
